$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37, shifting existing rows 37-64 down to 38-65.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Cells.Item(37, 1).Value = 10
$ws.Cells.Item(37, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37, 3).Value = "La Araucanía"
$ws.Cells.Item(37, 4).Value = 44957
$ws.Cells.Item(37, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100107
$ws.Cells.Item(37, 8).Value = "Otros"
$ws.Cells.Item(37, 9).Value = 100107011
$ws.Cells.Item(37, 10).Value = "Tuna"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 25
$ws.Cells.Item(37, 14).Value = 30000
$ws.Cells.Item(37, 15).Value = 30000
$ws.Cells.Item(37, 16).Value = 30000
$ws.Cells.Item(37, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(37, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(37, 19).Value = 1667
$ws.Cells.Item(37, 20).Value = 18
